# Update the born (revive/spawn) position for the "villageScene" (city)
# entry in the Scene config sheet.  RelivePos lives in column E; row 2 is
# the villageScene / city row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "20,0,-137"
